# Reverse the order of the comma-separated "Recorded By" entries in column G.
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Single-value cells (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.Contains(",")) {
            $parts = $text.Split(",")
            $reversed = ""
            for ($i = $parts.Length - 1; $i -ge 0; $i--) {
                $piece = $parts[$i].Trim()
                if ($reversed -ne "") {
                    $reversed = $reversed + ", "
                }
                $reversed = $reversed + $piece
            }
            $cell.Value2 = $reversed
        }
    }
}
